$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "MIE(7.953679543400868, 4.056257339577724, -9.373344930178234, 10.501397738292106)"
$ws.Range("C2").Value = "NIG(0.9902598653969242, 0.7481245965449379, 4.526136526342097, 5.708367079522397)"
$ws.Range("D2").Value = "JSU(-0.7941703336716979, 1.0124765125902626, 1.120878684712627, 2.410839405689891)"
$ws.Range("E2").Value = "NCT(3.445188548131043, 1.9152323435338388, -0.018231139591131984, 4.361234560057879)"
